# mastertimeseries.xlsx update:
#  - insert a new column before the existing "W" (Utahgasprice_dollperMCF) column
#  - populate the new column with header "emisCO2eq20_millnMg" and a CO2-eq
#    formula derived from column B (methane emissions), filled down for every
#    row that has a value in B (row 3 / year 2014 has no B value, so it is
#    left blank, matching the source data)
#  - add Seth Lyman's note on the new column explaining where the formula /
#    factor (82.5, from Logan Mitchell @ Utah Clean Energy) came from
#  - the two pre-existing header comments (gas price source, oil price
#    source) stay attached to their original cells, which slide one column
#    to the right (W1 -> X1, X1 -> Y1) along with the data

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- capture the text of the two existing header comments before we shift
# anything, so we can re-create them in the right place -----------------
$gasPriceComment = $ws.Range("W1").Comment.Text()
$oilPriceComment = $ws.Range("X1").Comment.Text()

# --- insert the new column in front of W (col 23); existing W/X data,
# formatting and comments (as far as values go) shift right one column --
$ws.Columns.Item(23).Insert()

# the comments on the old W1/X1 do not automatically travel with the
# insert in this host, so re-home them explicitly on the cells they now
# belong to (X1 / Y1) and drop whatever stale comment objects remain
if ($ws.Range("W1").Comment -ne $null) { $ws.Range("W1").Comment.Delete() }
if ($ws.Range("X1").Comment -ne $null) { $ws.Range("X1").Comment.Delete() }
if ($ws.Range("Y1").Comment -ne $null) { $ws.Range("Y1").Comment.Delete() }

$ws.Range("X1").AddComment() | Out-Null
$ws.Range("X1").Comment.Text($gasPriceComment)

$ws.Range("Y1").AddComment() | Out-Null
$ws.Range("Y1").Comment.Text($oilPriceComment)

# --- new column header / data -------------------------------------------
$ws.Range("W1").Value = "emisCO2eq20_millnMg"

$ws.Range("W2").Formula = "=(B2*24*365.25)*82.5/1000000"
# row 3 (year 2014) has no basinwide_ch4_emiss value in column B, so the
# new column is left blank there too, same as the rest of that row
$ws.Range("W4:W13").Formula = "=(B4*24*365.25)*82.5/1000000"

# --- new comment on the header, explaining the formula ------------------
$ws.Range("W1").AddComment() | Out-Null
$ws.Range("W1").Comment.Text("Seth Lyman:" + [char]10 + "=(B2*24*365.25)*82.5/1000000" + [char]10 + "from Logan Mitchell, utah clean energy")

# --- restore the "current selection" look of the sheet -------------------
$ws.Range("B2").Select()
$excel.ActiveWindow.FreezePanes = $false
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("W1").Select()
